# Auto-generated edit script applying the scheduled market-data refresh
# to the Adamantoise_Profits workbook (values for H..N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 636.6
$ws.Range("I33").Value = 394.4
$ws.Range("J33").Value = 1121
$ws.Range("K33").Value = 394.4
$ws.Range("L33").Value = 1121
$ws.Range("M33").Value = -165.4

# Row 62
$ws.Range("H62").Value = 10115
$ws.Range("I62").Value = 9000
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -8376
$ws.Range("N62").Value = -13593

# Row 64
$ws.Range("I64").Value = 8916.333000000001
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 8916.333000000001
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -8668.333000000001

# Row 65
$ws.Range("H65").Value = 10115
$ws.Range("I65").Value = 9000
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -67965

# Row 67
$ws.Range("I67").Value = 8916.333000000001
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 8916.333000000001
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -8058.333000000001

# Row 74
$ws.Range("H74").Value = 5383.091
$ws.Range("I74").Value = 6173.5713
$ws.Range("J74").Value = 3999.75
$ws.Range("K74").Value = 6173.5713
$ws.Range("L74").Value = 3999.75
$ws.Range("M74").Value = -5237.5713

# Row 77
$ws.Range("H77").Value = 5383.091
$ws.Range("I77").Value = 6173.5713
$ws.Range("J77").Value = 3999.75
$ws.Range("K77").Value = 30867.8565
$ws.Range("L77").Value = 19998.75
$ws.Range("M77").Value = -26187.8565

# Row 141
$ws.Range("H141").Value = 6184.3335
$ws.Range("I141").Value = 5180.2104
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 15540.6312
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -10360.6312

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3809.9375
$ws.Range("I61").Value = 3114.6667
$ws.Range("J61").Value = 4423.4116
$ws.Range("K61").Value = 3114.6667
$ws.Range("L61").Value = 4423.4116
$ws.Range("M61").Value = -2902.6667

# Row 63
$ws.Range("H63").Value = 3314.889
$ws.Range("I63").Value = 2197.3333
$ws.Range("J63").Value = 5550
$ws.Range("K63").Value = 2197.3333
$ws.Range("L63").Value = 5550
$ws.Range("M63").Value = -1511.3333
$ws.Range("N63").Value = -6922

# Row 66
$ws.Range("H66").Value = 3314.889
$ws.Range("I66").Value = 2197.3333
$ws.Range("J66").Value = 5550
$ws.Range("K66").Value = 10986.6665
$ws.Range("L66").Value = 27750
$ws.Range("M66").Value = -7554.666499999999
$ws.Range("N66").Value = -34614

# Row 122
$ws.Range("H122").Value = 5792.1816
$ws.Range("I122").Value = 5321.4287
$ws.Range("J122").Value = 6011.8667
$ws.Range("K122").Value = 15964.2861
$ws.Range("L122").Value = 18035.6001
$ws.Range("M122").Value = -13514.2861
$ws.Range("N122").Value = -22935.6001

# Row 132
$ws.Range("H132").Value = 305886.1
$ws.Range("I132").Value = 305886.1
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 917658.2999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -915128.2999999999
$ws.Range("N132").ClearContents()

# Row 136
$ws.Range("H136").Value = 3809.9375
$ws.Range("I136").Value = 3114.6667
$ws.Range("J136").Value = 4423.4116
$ws.Range("K136").Value = 9344.000100000001
$ws.Range("L136").Value = 13270.2348
$ws.Range("M136").Value = -6794.000100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 45104.082
$ws.Range("I20").Value = 103111
$ws.Range("J20").Value = 3670.5715
$ws.Range("K20").Value = 103111
$ws.Range("L20").Value = 3670.5715
$ws.Range("M20").Value = -102864
$ws.Range("N20").Value = -4164.5715

# Row 134
$ws.Range("I134").Value = 6064244.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18192733.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -18190198.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 250249
$ws.Range("I4").Value = 250249
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 250249
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -250137

# Row 31
$ws.Range("H31").Value = 3521.3547
$ws.Range("I31").Value = 1338.9375
$ws.Range("J31").Value = 5849.2666
$ws.Range("K31").Value = 1338.9375
$ws.Range("L31").Value = 5849.2666
$ws.Range("M31").Value = -1043.9375
$ws.Range("N31").Value = -6439.2666

# Row 34
$ws.Range("H34").Value = 3521.3547
$ws.Range("I34").Value = 1338.9375
$ws.Range("J34").Value = 5849.2666
$ws.Range("K34").Value = 1338.9375
$ws.Range("L34").Value = 5849.2666
$ws.Range("M34").Value = -1136.9375
$ws.Range("N34").Value = -6253.2666

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()

# Row 94
$ws.Range("H94").Value = 1349.1177
$ws.Range("I94").Value = 864.3333
$ws.Range("J94").Value = 1453
$ws.Range("K94").Value = 864.3333
$ws.Range("L94").Value = 1453
$ws.Range("M94").Value = -413.3333
$ws.Range("N94").Value = -2355

# Row 114
$ws.Range("H114").Value = 18250
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 18250
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 18250
$ws.Range("N114").Value = -26928

# Row 134
$ws.Range("H134").Value = 2232.5
$ws.Range("I134").Value = 2279.2
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 6837.599999999999
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -4302.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 127
$ws.Range("H127").Value = 2445
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2445
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 7335
$ws.Range("N127").Value = -17255

# Row 139
$ws.Range("H139").Value = 4065.25
$ws.Range("I139").Value = 3646
$ws.Range("J139").Value = 7000
$ws.Range("K139").Value = 10938
$ws.Range("L139").Value = 21000
$ws.Range("M139").Value = -5798

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2559
$ws.Range("I80").Value = 2439.3333
$ws.Range("J80").Value = 2918
$ws.Range("K80").Value = 2439.3333
$ws.Range("L80").Value = 2918
$ws.Range("M80").Value = -1441.3333
$ws.Range("N80").Value = -4914

# Row 83
$ws.Range("H83").Value = 2559
$ws.Range("I83").Value = 2439.3333
$ws.Range("J83").Value = 2918
$ws.Range("K83").Value = 12196.6665
$ws.Range("L83").Value = 14590
$ws.Range("M83").Value = -7204.666499999999
$ws.Range("N83").Value = -24574

# Row 122
$ws.Range("H122").Value = 5474.75
$ws.Range("I122").Value = 1949.5
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 5848.5
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -3398.5

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 82
$ws.Range("H82").Value = 3313.8235
$ws.Range("I82").Value = 3857.4285
$ws.Range("J82").Value = 2933.3
$ws.Range("K82").Value = 3857.4285
$ws.Range("L82").Value = 2933.3
$ws.Range("M82").Value = -3496.4285

# Row 85
$ws.Range("H85").Value = 3313.8235
$ws.Range("I85").Value = 3857.4285
$ws.Range("J85").Value = 2933.3
$ws.Range("K85").Value = 3857.4285
$ws.Range("L85").Value = 2933.3
$ws.Range("M85").Value = -2609.4285

# Row 122
$ws.Range("H122").Value = 20576.166
$ws.Range("I122").Value = 24000.285
$ws.Range("J122").Value = 15782.4
$ws.Range("K122").Value = 72000.855
$ws.Range("L122").Value = 47347.2
$ws.Range("M122").Value = -69550.855

# Row 136
$ws.Range("H136").Value = 7125.0835
$ws.Range("I136").Value = 2072
$ws.Range("J136").Value = 14199.4
$ws.Range("K136").Value = 6216
$ws.Range("L136").Value = 42598.2
$ws.Range("M136").Value = -3666
$ws.Range("N136").Value = -47698.2

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 5081.3335
$ws.Range("I54").Value = 5081.3335
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5081.3335
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4561.3335

# Row 132
$ws.Range("H132").Value = 33334.637
$ws.Range("I132").Value = 38330.32
$ws.Range("J132").Value = 5358.8
$ws.Range("K132").Value = 114990.96
$ws.Range("L132").Value = 16076.4
$ws.Range("M132").Value = -112460.96

# Row 136
$ws.Range("H136").Value = 2294.25
$ws.Range("I136").Value = 1433.2632
$ws.Range("J136").Value = 3552.6155
$ws.Range("K136").Value = 4299.7896
$ws.Range("L136").Value = 10657.8465
$ws.Range("M136").Value = -1749.7896
$ws.Range("N136").Value = -15757.8465
